$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Heap (Basics)" section (old rows 37 + 39) needs to move down
# to make room for a brand-new "Sliding Window (Advance)" section
# (new rows 37 + 39). Insert 5 blank rows before the old row 37 so
# that old row 37 -> 42 and old row 39 -> 44, matching the target
# layout (rows 40/41/43 stay empty, exactly like row 38 already is).
# ------------------------------------------------------------------
$ws.Rows.Item(37).Resize(5).Insert()

# Row 36: blank filler row (matches the style already used one row
# above it, exactly like every other blank spacer row in this sheet).
$ws.Range("B36").Value = ""
$ws.Range("B36").Font.ThemeColor = $ws.Range("B35").Font.ThemeColor
$ws.Range("B36").Font.Size = $ws.Range("B35").Font.Size
$ws.Range("E36").NumberFormat = $ws.Range("E35").NumberFormat

# Row 37: new section header "Sliding Window (Advance)" - same look
# as the other plain section headers (B14/B22/B27/B32), i.e. themed
# fill + bold-ish header font, taller row.
$ws.Range("B37").Value = "Sliding Window (Advance)"
$ws.Range("B37").Font.ThemeColor = $ws.Range("B32").Font.ThemeColor
$ws.Range("B37").Font.Size = $ws.Range("B32").Font.Size
$ws.Range("B37").Interior.ThemeColor = $ws.Range("B32").Interior.ThemeColor
$ws.Rows.Item(37).RowHeight = 39

# Row 39: new data row for question 424.
$ws.Range("A39").Value = 424
$ws.Range("B39").Value = "Longest Repeating Character Replacement"
$ws.Range("B39").WrapText = $true
$ws.Range("C39").Value = "Medium"
$ws.Range("D39").Value = "String,sliding window,hashmap"
$ws.Range("E39").Value = 45716
$ws.Range("E39").NumberFormat = $ws.Range("E35").NumberFormat
$ws.Rows.Item(39).RowHeight = 17

# Selection ends up on E40, matching the author's last-touched cell.
$ws.Range("E40").Select()

Write-Output "done"
